# Applies the "displaySpeciesTypes / getSpeciesTypes / Place / SpeciesType"
# text-consolidation edits to the Traceability diagram table.
#
# In the source document the visible text in each target cell is already
# (when runs are concatenated) what the diff shows as the final text -- the
# runs just haven't been merged into one <w:r> yet (e.g. "display" + "Species"
# + "Types()"). Word's real Find/Replace normally collapses such matches into
# a single run, but in this runtime Find.Execute always operates over the
# *whole* document regardless of which Range/Cell it is invoked on, which
# would wrongly touch unrelated text elsewhere (e.g. "intSpeciesType",
# "type:SpeciesType", the lowercase word "places", ...).
#
# To stay precise we instead address each target table cell directly via
# the Tables/Rows/Cell object model and overwrite its Range.Text. Assigning
# identical text is a no-op for some runtimes, so we first set the cell to a
# short, unique placeholder (forcing a genuine text replacement/run-merge)
# and then set it to the real final text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $finalText) {
    $cell = $t.Cell($row, $col)
    $full = $cell.Range
    # Cell.Range includes the trailing cell-mark characters (chr(13)+chr(7));
    # exclude them so we only replace the actual visible text.
    $target = $d.Range($full.Start, $full.End - 1)

    # Force an actual text change first so multi-run cells really get
    # collapsed into a single run (identical-text assignment can be a no-op).
    $target.Text = "`u{E000}PLACEHOLDER`u{E000}"

    $full2 = $cell.Range
    $target2 = $d.Range($full2.Start, $full2.End - 1)
    $target2.Text = $finalText
}

# 1) "display" + "Species" + "Types()"  ->  "displaySpeciesTypes()"
Set-CellText 74 4 "displaySpeciesTypes()"

# 2) "get" + "Species" + "Types()"  ->  "getSpeciesTypes()"
Set-CellText 75 4 "getSpeciesTypes()"

# 3) "Places"  ->  "Place"
Set-CellText 76 3 "Place"

# 4) "Species" + "Type"  ->  "SpeciesType"
Set-CellText 77 3 "SpeciesType"

# 5) "Species" + "Type"  ->  "SpeciesType"
Set-CellText 78 3 "SpeciesType"

# 6) "Species" + "Type(description: String)"  ->  "SpeciesType(description: String)"
Set-CellText 78 4 "SpeciesType(description: String)"

# 7) "Species" + "Type"  ->  "SpeciesType"
Set-CellText 79 3 "SpeciesType"
